# run_program() name output files with ordinal number
# Insert a new first column ("lp" = ordinal number in Polish) before the
# existing "NR faktury" / "nazwa pliku" columns, and number the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:B to B:C, inserting a new blank column A.
$ws.Columns("A").Insert()

# Header for the new ordinal-number column.
$ws.Range("A1").Value2 = "lp"

# Sequential ordinal numbers for the data rows (rows 2-7).
$ws.Range("A2").Value2 = 1
$ws.Range("A3").Value2 = 2
$ws.Range("A4").Value2 = 3
$ws.Range("A5").Value2 = 4
$ws.Range("A6").Value2 = 5
$ws.Range("A7").Value2 = 6

# Row 8 (a differently-shaped record) is left without an ordinal number.

# Match the resulting active selection.
$ws.Range("A8").Select()
